# Insert a new weekly price record at row 372 (Zanahoria, Terminal La Palmera
# de La Serena). This pushes the former rows 372-492 down to 373-493.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(372).EntireRow.Insert()

$ws.Cells.Item(372, 1).Value  = 8
$ws.Cells.Item(372, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(372, 3).Value  = "Coquimbo"
$ws.Cells.Item(372, 4).Value  = 44988
$ws.Cells.Item(372, 5).Value  = 4
$ws.Cells.Item(372, 6).Value  = 100114013
$ws.Cells.Item(372, 7).Value  = "Zanahoria"
$ws.Cells.Item(372, 8).Value  = "Sin especificar"
$ws.Cells.Item(372, 9).Value  = "Primera"
$ws.Cells.Item(372, 10).Value = 400
$ws.Cells.Item(372, 11).Value = 5500
$ws.Cells.Item(372, 12).Value = 6000
$ws.Cells.Item(372, 13).Value = 5750
$ws.Cells.Item(372, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(372, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(372, 16).Value = 288
$ws.Cells.Item(372, 17).Value = 20
$ws.Cells.Item(372, 18).Value = "Hortaliza"
